# Apply updated dSF (column F) values on Sheet1, per commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value  = -4
$ws.Range("F5").Value  = -7
$ws.Range("F6").Value  = -3
$ws.Range("F7").Value  = -3
$ws.Range("F8").Value  = 1
$ws.Range("F10").Value = -1
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -7
$ws.Range("F17").Value = -3
$ws.Range("F19").Value = -8
